$d = $word.ActiveDocument

# ===========================================================================
# Background: the "Question 15" heading paragraph ("15. What is the purpose
# of the GROUP BY clause, and how is it used? ") is currently split across
# three runs ("15" / "." / " What is the purpose of the GROUP BY clause,
# and how is it used? ") that all share identical bold / size-28
# formatting. The very next paragraph is its (unformatted) answer.
#
# Target state:
#   - The Question-15 heading collapses into a single bold run holding the
#     full "15. What is ... used? " text (no visible change, just a run
#     merge).
#   - A brand-new bold heading paragraph for "16. What is the WHERE clause
#     used for, and how is it used to filter data? " is inserted right
#     after the GROUP BY answer paragraph.
#   - A brand-new plain-text answer paragraph about the WHERE clause is
#     inserted right after that new heading.
# ===========================================================================

# ---------------------------------------------------------------------------
# Locate the Question-15 heading paragraph by its (unique) text rather than
# a hard-coded index.
# ---------------------------------------------------------------------------
$q15Index = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*What is the purpose of the GROUP BY clause*") {
        $q15Index = $i
        break
    }
}

# ---------------------------------------------------------------------------
# 1) Merge the Question-15 heading's three runs into a single bold run.
#    Strategy: capture its plain text, delete it, retype it as plain text,
#    then apply Bold/Size to the freshly inserted span. Re-typing (instead
#    of re-using the old runs) is what lets Word collapse the three
#    original runs down into exactly one.
# ---------------------------------------------------------------------------
$q15 = $d.Paragraphs($q15Index)
$q15Range = $d.Range($q15.Range.Start, $q15.Range.End - 1)
$q15Text = $q15Range.Text
$q15Start = $q15Range.Start
$q15Range.Delete()

$q15Insert = $d.Range($q15Start, $q15Start)
$q15Insert.InsertAfter($q15Text)

$q15Formatted = $d.Range($q15Start, $q15Start + $q15Text.Length)
$q15Formatted.Font.Bold = 1
$q15Formatted.Font.Size = 14

# ---------------------------------------------------------------------------
# 2) Insert the new Question-16 heading paragraph and its answer paragraph
#    right after the GROUP BY answer paragraph (the paragraph immediately
#    following the Question-15 heading).
#
#    Important: insert ALL the new plain text first (paragraph breaks
#    included) *before* touching any Bold/Size formatting. If formatting is
#    applied too early, later InsertAfter calls / InsertParagraphAfter
#    splits pick up that formatting as ambient context and it leaks onto
#    text that must stay in the default (un-bolded) style.
# ---------------------------------------------------------------------------
$groupByAnswerIndex = $q15Index + 1
$groupByAnswer = $d.Paragraphs($groupByAnswerIndex)
$endOfAnswer = $d.Range($groupByAnswer.Range.Start, $groupByAnswer.Range.End - 1)
$endOfAnswer.Collapse(0)
$endOfAnswer.InsertParagraphAfter()

# New empty heading paragraph immediately follows the GROUP BY answer.
$q16Index = $groupByAnswerIndex + 1
$q16Heading = $d.Paragraphs($q16Index)
$q16HeadingInsert = $d.Range($q16Heading.Range.Start, $q16Heading.Range.Start)
$q16HeadingText = "16. What is the WHERE clause used for, and how is it used to filter data? "
$q16HeadingInsert.InsertAfter($q16HeadingText)

# Break off the answer paragraph right after the heading (still plain text).
$q16HeadingAfterInsert = $d.Paragraphs($q16Index)
$endOfHeading = $d.Range($q16HeadingAfterInsert.Range.Start, $q16HeadingAfterInsert.Range.End - 1)
$endOfHeading.Collapse(0)
$endOfHeading.InsertParagraphAfter()

$q16AnswerIndex = $q16Index + 1
$q16Answer = $d.Paragraphs($q16AnswerIndex)
$q16AnswerInsert = $d.Range($q16Answer.Range.Start, $q16Answer.Range.Start)
$q16AnswerText = "In SQL, the WHERE clause is used to filter rows from a table or result set according to predetermined criteria. It enables us to pick only the rows that satisfy particular requirements or follow a pattern. A key element of SQL queries, the WHERE clause is frequently used for data retrieval and manipulation."
$q16AnswerInsert.InsertAfter($q16AnswerText)

# Now apply bold/size-28 formatting to just the heading's bold portion,
# i.e. everything except the trailing space (which stays unformatted, just
# like the rest of the document's body text).
$q16HeadingFinal = $d.Paragraphs($q16Index)
$boldLength = $q16HeadingText.Length - 1
$q16HeadingBoldRange = $d.Range($q16HeadingFinal.Range.Start, $q16HeadingFinal.Range.Start + $boldLength)
$q16HeadingBoldRange.Font.Bold = 1
$q16HeadingBoldRange.Font.Size = 14

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
